# Re-ran the sweep over the full map (all 17 configs). Updates the
# existing result rows (cols G:L -- count/mean/var/skew/avoid-flag/dist)
# with the new measurements and appends the 3 extra rows the wider run
# produced. Avoid-policy cascade (col K/L) still looks shaky on a few rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @{row=1; col=8; val=30.8169106},
  @{row=1; col=9; val=0.0019058552724922073},
  @{row=1; col=10; val=0.000018392708652162348},
  @{row=2; col=7; val=200.0},
  @{row=2; col=8; val=81.300641},
  @{row=2; col=9; val=0.0011763871526619951},
  @{row=2; col=10; val=0.000013793771315107617},
  @{row=2; col=12; val=27.0},
  @{row=3; col=8; val=75.6373166},
  @{row=3; col=9; val=0.0011763871526619951},
  @{row=3; col=10; val=0.000013561149916247957},
  @{row=4; col=7; val=186.0},
  @{row=4; col=8; val=165.3663291},
  @{row=4; col=9; val=0.0011763871526619951},
  @{row=4; col=10; val=0.000013957391818017622},
  @{row=5; col=8; val=85.326889},
  @{row=5; col=9; val=0.001390708605742308},
  @{row=5; col=10; val=0.000020304716180701956},
  @{row=5; col=12; val=25.0},
  @{row=6; col=7; val=428.0},
  @{row=6; col=8; val=257.6468557},
  @{row=6; col=9; val=0.0024080313284957633},
  @{row=6; col=10; val=0.000013723826764045608},
  @{row=7; col=7; val=487.0},
  @{row=7; col=8; val=204.9694213},
  @{row=7; col=9; val=0.002946943668191704},
  @{row=7; col=10; val=0.000016598221994684326},
  @{row=7; col=12; val=17.0},
  @{row=8; col=7; val=116.0},
  @{row=8; col=8; val=63.4202869},
  @{row=8; col=9; val=0.0003260316997157897},
  @{row=8; col=10; val=-0.0000027236118733130056},
  @{row=8; col=12; val=23.0},
  @{row=9; col=7; val=165.0},
  @{row=9; col=8; val=94.326619},
  @{row=9; col=9; val=0.0072346880459222},
  @{row=9; col=10; val=0.00004441840052835521},
  @{row=9; col=11; val=3.0},
  @{row=9; col=12; val=48.0},
  @{row=10; col=7; val=179.0},
  @{row=10; col=8; val=146.1809193},
  @{row=10; col=9; val=0.0003782526172055878},
  @{row=10; col=10; val=-0.000027389651869448497},
  @{row=10; col=11; val=2.0},
  @{row=10; col=12; val=57.0},
  @{row=11; col=7; val=204.0},
  @{row=11; col=8; val=115.3778066},
  @{row=11; col=9; val=0.0009027301200268401},
  @{row=11; col=10; val=0.000013570447126830444},
  @{row=12; col=7; val=180.0},
  @{row=12; col=8; val=116.0721434},
  @{row=12; col=9; val=0.001524599347155542},
  @{row=12; col=10; val=-0.000028272627787815734},
  @{row=12; col=11; val=2.0},
  @{row=12; col=12; val=57.0},
  @{row=13; col=7; val=176.0},
  @{row=13; col=8; val=116.3119104},
  @{row=13; col=9; val=0.0019422957481607384},
  @{row=13; col=10; val=-0.000012999783357839103},
  @{row=13; col=11; val=2.0},
  @{row=13; col=12; val=57.0},
  @{row=14; col=7; val=150.0},
  @{row=14; col=8; val=27.7283103},
  @{row=14; col=9; val=0.001390708605742308},
  @{row=14; col=10; val=0.00001994476897638468},
  @{row=14; col=11; val=0.0},
  @{row=14; col=12; val=0.0}
)
foreach ($u in $updates) {
  $ws.Cells.Item($u.row, $u.col).Value = $u.val
}

$newRows = @(
  @(10.0, 1000.0, 5.0, 2.5, 100.0, 0.95, 153.0, 64.4054013, 0.001390708605742308, 0.000020256190713535646, 0.0, 0.0),
  @(10.0, 1000.0, 5.0, 2.5, 100.0, 0.95, 165.0, 94.3117423, 0.001190702130974275, 0.000015755517712114045, 0.0, 0.0),
  @(10.0, 1000.0, 5.0, 2.5, 100.0, 0.95, 182.0, 71.5354921, 0.0023318558688201207, -0.00007349815196472742, 1.0, 30.0)
)
$startRow = 15
for ($i = 0; $i -lt $newRows.Count; $i++) {
  $rowVals = $newRows[$i]
  $r = $startRow + $i
  for ($c = 1; $c -le $rowVals.Count; $c++) {
    $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
  }
}
